# "Added last minute updates"
#
# The first paragraph in the document contains a docx4j merge-field
# placeholder. This edit:
#   1. Updates the placeholder ID text
#      (**ID__AFFARS_pgi_5309_topic_7__ID** -> **ID__AFFARS_AFMC_PGI_5309_405__ID**)
#      and drops the now-redundant trailing-space run that used to follow it.
#   2. Adds an (empty/no-line) paragraph border with 5-twip spacing on all
#      four sides to that paragraph's pPr.
#   3. Widens the paragraph's left indent from 120 to 225 twips.

$d = $word.ActiveDocument

$p1 = $d.Paragraphs.Item(1)

# --- Paragraph border (space-only, no visible line) -----------------------
# Setting DistanceFromXxx on the Borders collection writes <w:top w:space="5"/>
# style entries (no w:val/w:sz/w:color) without touching the run contents,
# which keeps the existing run formatting byte-for-byte intact.
$borders = $p1.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# --- Left indent: 120 -> 225 twips (twips / 20 = points) -------------------
$p1.Format.LeftIndent = 11.25

# --- Replace the placeholder text and swallow the trailing space run -------
# Matching through the trailing space merges it into the first run and
# removes the now-empty second run entirely.
$d.Content.Find.Execute("**ID__AFFARS_pgi_5309_topic_7__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5309_405__ID**", 2)
